# Fill in the "Hours Left" burndown values for days 1-3 so the working
# diction persists between sessions instead of resetting to 0.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 13
$ws.Range("B4").Value = 12
$ws.Range("B5").Value = 11

# Restore the active selection to B6 (was previously left at H18).
$ws.Range("B6").Select()
